# Add a new Job Posting row (JD_008 - Senior System Engineer) to the
# LinkedIn job postings sheet, as the last row right after JD_007.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row below the existing data (row 9, since rows 1-8
# are already populated: 1 header + 7 job postings).
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$ws.Cells.Item($newRow, 1).Value = "JD_008"
$ws.Cells.Item($newRow, 2).Value = "Senior System Engineer"
$ws.Cells.Item($newRow, 3).Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Cells.Item($newRow, 4).Value = 2
$ws.Cells.Item($newRow, 5).Value = 4

# Reset the auto-grown row height back to the sheet default (matches the
# other data rows, which have no explicit/custom row height).
$ws.Rows($newRow).EntireRow.AutoFit()
